$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: product changes from GPT-0532 (OKI TONER...) to GPI-0224 (CANON INK...)
$ws.Range("I2").Value = "GPI-0224"
$ws.Range("J2").Value = "CANON INK CLI-551XL GRAY ΣΥΜΒΑΤΟ 13ml"

$ws.Range("K2").Formula = "=""23"""
$ws.Range("K2").Copy()
$ws.Range("K2").PasteSpecial(-4163)

$ws.Range("L2").Formula = "=""780€"""
$ws.Range("L2").Copy()
$ws.Range("L2").PasteSpecial(-4163)

# Row 3: product changes from GPI-0134 (HP INK No 88XL...) to GPI-0023 (EPSON INK No 26XL...)
$ws.Range("I3").Value = "GPI-0023"
$ws.Range("J3").Value = "EPSON INK No 26XL - T2634XL YELLOW ΣΥΜΒΑΤΟ 10ml"

$ws.Range("K3").Formula = "=""16"""
$ws.Range("K3").Copy()
$ws.Range("K3").PasteSpecial(-4163)

# Row 4: product changes from GPI-0023 (EPSON INK No 26XL...) to GPT-0070 (EPSON TONER EPL6200L...)
$ws.Range("I4").Value = "GPT-0070"
$ws.Range("J4").Value = "EPSON TONER EPL6200L BLACK ΣΥΜΒΑΤΟ 6000 ΣΕΛΙΔΕΣ"

$ws.Range("K4").Formula = "=""120"""
$ws.Range("K4").Copy()
$ws.Range("K4").PasteSpecial(-4163)

# Row 5 is removed entirely, shrinking the used range to A1:L4
$ws.Rows.Item(5).Delete()
